# Completed DeleteCustomer, CreateSO, LogActivity, LogCommunication, ActionPerform by Nitin.
#
# The test-automation log sheet ("Customers_details") picks up a new
# "last-run" timestamp in C2 plus nine appended columns (AB:AJ) holding
# Invoice / Payment / Expense data captured while the new scenarios ran.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- refreshed run timestamp -------------------------------------------------
$ws.Range("C2").Value = "Quick 2019/11/26 18:58:19"

# --- new header row (AB1:AJ1) --------------------------------------------
$ws.Range("AB1").Value = "Invoice Number"
$ws.Range("AC1").Value = "Payment mode"
$ws.Range("AD1").Value = "Payment Mode Details"
$ws.Range("AE1").Value = "Payment Amount"
$ws.Range("AF1").Value = "Address Name"
$ws.Range("AG1").Value = "Country"
$ws.Range("AH1").Value = "Expense Type"
$ws.Range("AI1").Value = "Expense Amount"
$ws.Range("AJ1").Value = "Expense Number"

# --- new data row (AB2:AI2) -----------------------------------------------
$ws.Range("AB2").Value = "EXP-526-261119-8"
$ws.Range("AB2").Style = "Normal"

# "300" must land as literal text (matches the source file's <t>300</t>
# shared-string entry), not as a number -- a plain .Value assignment gets
# auto-coerced to numeric. Route it through a text-formula + paste-values
# round-trip so it sticks as text without picking up a new number format
# / style entry.
$ws.Range("AC2").Formula = "=""300"""
$ws.Range("AC2").Copy()
$ws.Range("AC2").PasteSpecial(-4163)
$ws.Range("AC2").Style = "Normal"

$ws.Range("AD2").Value = "Online Bank"
$ws.Range("AE2").Value = 199
$ws.Range("AF2").Value = "Local"
$ws.Range("AG2").Value = "India"
$ws.Range("AH2").Value = "Shipping Expense"
$ws.Range("AI2").Value = 100

# --- selection / view state -------------------------------------------------
$ws.Select()
$ws.Range("AE8").Select()
$excel.ActiveWindow.Width = 14690
